$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 15
$ws.Cells.Item(15, 1).Value = 111356613
$ws.Cells.Item(15, 2).Value = 88819
$ws.Cells.Item(15, 4).Value = 'LC'
$ws.Cells.Item(15, 5).Value = 5685
$ws.Cells.Item(15, 6).Value = 'Gullgröppa'
$ws.Cells.Item(15, 7).Value = 'Pseudomerulius aureus'
$ws.Cells.Item(15, 8).Value = '(Fr.) Jülich'
$ws.Cells.Item(15, 9).ClearContents()
$ws.Cells.Item(15, 17).Value = 491932.6359417734
$ws.Cells.Item(15, 18).Value = 6785481.417835217
$ws.Cells.Item(15, 19).Value = 10

# Row 16
$ws.Cells.Item(16, 1).Value = 111356614
$ws.Cells.Item(16, 2).Value = 96348
$ws.Cells.Item(16, 4).Value = 'VU'
$ws.Cells.Item(16, 5).Value = 220787
$ws.Cells.Item(16, 6).Value = 'Knärot'
$ws.Cells.Item(16, 7).Value = 'Goodyera repens'
$ws.Cells.Item(16, 8).Value = '(L.) R. Br.'
$ws.Cells.Item(16, 9).Value = '15'
$ws.Cells.Item(16, 17).Value = 491909.0202035823
$ws.Cells.Item(16, 18).Value = 6785498.341940038
$ws.Cells.Item(16, 19).Value = 25

# Row 19
$ws.Cells.Item(19, 1).Value = 111491680
$ws.Cells.Item(19, 2).Value = 56414
$ws.Cells.Item(19, 4).Value = 'NT'
$ws.Cells.Item(19, 5).Value = 100049
$ws.Cells.Item(19, 6).Value = 'Spillkråka'
$ws.Cells.Item(19, 7).Value = 'Dryocopus martius'
$ws.Cells.Item(19, 8).Value = '(Linnaeus, 1758)'
$ws.Cells.Item(19, 10).ClearContents()
$ws.Cells.Item(19, 12).Value = ''
$ws.Cells.Item(19, 13).Value = 'äldre spår'
$ws.Cells.Item(19, 17).Value = 491929.8523854768
$ws.Cells.Item(19, 18).Value = 6785530.587422797
$ws.Cells.Item(19, 29).ClearContents()
$ws.Cells.Item(19, 32).ClearContents()

# Row 20
$ws.Cells.Item(20, 1).Value = 111491685
$ws.Cells.Item(20, 2).Value = 88819
$ws.Cells.Item(20, 4).Value = 'LC'
$ws.Cells.Item(20, 5).Value = 5685
$ws.Cells.Item(20, 6).Value = 'Gullgröppa'
$ws.Cells.Item(20, 7).Value = 'Pseudomerulius aureus'
$ws.Cells.Item(20, 8).Value = '(Fr.) Jülich'
$ws.Cells.Item(20, 10).Value = ''
$ws.Cells.Item(20, 12).ClearContents()
$ws.Cells.Item(20, 13).ClearContents()
$ws.Cells.Item(20, 17).Value = 491909.4940688942
$ws.Cells.Item(20, 18).Value = 6785494.484901348
$ws.Cells.Item(20, 32).Value = ''

# Row 21
$ws.Cells.Item(21, 1).Value = 111491639
$ws.Cells.Item(21, 2).Value = 8377
$ws.Cells.Item(21, 5).Value = 106545
$ws.Cells.Item(21, 6).Value = 'Mindre märgborre'
$ws.Cells.Item(21, 7).Value = 'Tomicus minor'
$ws.Cells.Item(21, 8).Value = '(Hartig, 1834)'
$ws.Cells.Item(21, 12).Value = ''
$ws.Cells.Item(21, 13).Value = 'äldre gnagspår'
$ws.Cells.Item(21, 17).Value = 491993.9996831641
$ws.Cells.Item(21, 18).Value = 6785505.377163783
$ws.Cells.Item(21, 19).Value = 100
$ws.Cells.Item(21, 29).Value = 'Spridd och riklig i området'

# Row 23
$ws.Cells.Item(23, 1).Value = 111491662
$ws.Cells.Item(23, 2).Value = 77515
$ws.Cells.Item(23, 4).Value = 'NT'
$ws.Cells.Item(23, 5).Value = 6425
$ws.Cells.Item(23, 6).Value = 'Garnlav'
$ws.Cells.Item(23, 7).Value = 'Alectoria sarmentosa'
$ws.Cells.Item(23, 8).Value = '(Ach.) Ach.'
$ws.Cells.Item(23, 17).Value = 491969.9043757546
$ws.Cells.Item(23, 18).Value = 6785523.747805699

# Row 24
$ws.Cells.Item(24, 1).Value = 111491649
$ws.Cells.Item(24, 2).Value = 90666
$ws.Cells.Item(24, 4).Value = 'LC'
$ws.Cells.Item(24, 5).Value = 4364
$ws.Cells.Item(24, 6).Value = 'Dropptaggsvamp'
$ws.Cells.Item(24, 7).Value = 'Hydnellum ferrugineum'
$ws.Cells.Item(24, 8).Value = '(Fr.:Fr.) P. Karst.'
$ws.Cells.Item(24, 17).Value = 491979.6153062462
$ws.Cells.Item(24, 18).Value = 6785548.307010972
$ws.Cells.Item(24, 29).Value = 'mycel'

# Row 26
$ws.Cells.Item(26, 1).Value = 111491681
$ws.Cells.Item(26, 2).Value = 90666
$ws.Cells.Item(26, 5).Value = 4364
$ws.Cells.Item(26, 6).Value = 'Dropptaggsvamp'
$ws.Cells.Item(26, 7).Value = 'Hydnellum ferrugineum'
$ws.Cells.Item(26, 8).Value = '(Fr.:Fr.) P. Karst.'
$ws.Cells.Item(26, 12).ClearContents()
$ws.Cells.Item(26, 13).ClearContents()
$ws.Cells.Item(26, 17).Value = 491929.8523854768
$ws.Cells.Item(26, 18).Value = 6785530.587422797
$ws.Cells.Item(26, 19).Value = 5
$ws.Cells.Item(26, 29).Value = 'äldre fruktkreopp'

# Row 28
$ws.Cells.Item(28, 1).Value = 111612726
$ws.Cells.Item(28, 2).Value = 90168
$ws.Cells.Item(28, 4).Value = 'VU'
$ws.Cells.Item(28, 5).Value = 717
$ws.Cells.Item(28, 6).Value = 'Borsttagging'
$ws.Cells.Item(28, 7).Value = 'Gloiodon strigosus'
$ws.Cells.Item(28, 8).Value = '(Schwein. : Fr.) P. Karst.'
$ws.Cells.Item(28, 10).Value = ''
$ws.Cells.Item(28, 12).ClearContents()
$ws.Cells.Item(28, 13).ClearContents()
$ws.Cells.Item(28, 32).Value = ''

# Row 29
$ws.Cells.Item(29, 1).Value = 111612736
$ws.Cells.Item(29, 2).Value = 56398
$ws.Cells.Item(29, 5).Value = 100109
$ws.Cells.Item(29, 6).Value = 'Tretåig hackspett'
$ws.Cells.Item(29, 7).Value = 'Picoides tridactylus'

# Row 30
$ws.Cells.Item(30, 1).Value = 111612738
$ws.Cells.Item(30, 2).Value = 56414
$ws.Cells.Item(30, 4).Value = 'NT'
$ws.Cells.Item(30, 5).Value = 100049
$ws.Cells.Item(30, 6).Value = 'Spillkråka'
$ws.Cells.Item(30, 7).Value = 'Dryocopus martius'
$ws.Cells.Item(30, 8).Value = '(Linnaeus, 1758)'
$ws.Cells.Item(30, 10).ClearContents()
$ws.Cells.Item(30, 12).Value = ''
$ws.Cells.Item(30, 13).Value = 'äldre spår'
$ws.Cells.Item(30, 32).ClearContents()

# Row 32
$ws.Cells.Item(32, 1).Value = 111682652
$ws.Cells.Item(32, 2).Value = 90709
$ws.Cells.Item(32, 5).Value = 5448
$ws.Cells.Item(32, 6).Value = 'Svartvit taggsvamp'
$ws.Cells.Item(32, 7).Value = 'Phellodon connatus'
$ws.Cells.Item(32, 8).Value = '(Schultz) nom.prov'

# Row 33
$ws.Cells.Item(33, 1).Value = 111682769
$ws.Cells.Item(33, 2).Value = 89980
$ws.Cells.Item(33, 4).Value = 'VU'
$ws.Cells.Item(33, 5).Value = 1179
$ws.Cells.Item(33, 6).Value = 'Gräddticka'
$ws.Cells.Item(33, 7).Value = 'Perenniporia subacida'
$ws.Cells.Item(33, 8).Value = '(Peck) Donk'
$ws.Cells.Item(33, 17).Value = 491952.3910193561
$ws.Cells.Item(33, 18).Value = 6785464.984647369
$ws.Cells.Item(33, 29).ClearContents()
$ws.Cells.Item(33, 36).Value = 'gran'
$ws.Cells.Item(33, 37).Value = 'Picea abies'
$ws.Cells.Item(33, 41).Value = 'Picea abies'

# Row 34
$ws.Cells.Item(34, 1).Value = 111682665
$ws.Cells.Item(34, 2).Value = 90682
$ws.Cells.Item(34, 4).Value = 'NT'
$ws.Cells.Item(34, 5).Value = 2059
$ws.Cells.Item(34, 6).Value = 'Skrovlig taggsvamp'
$ws.Cells.Item(34, 7).Value = 'Hydnellum scabrosum'
$ws.Cells.Item(34, 8).Value = '(Fr.) E.Larss., K.H.Larss. & Kõljalg'
$ws.Cells.Item(34, 17).Value = 492024.0709204427
$ws.Cells.Item(34, 18).Value = 6785567.485207787
$ws.Cells.Item(34, 29).Value = 'Längs stigen/traktorspåret strax utanför gränsmarkeringen som syns på träden'
$ws.Cells.Item(34, 36).ClearContents()
$ws.Cells.Item(34, 37).ClearContents()
$ws.Cells.Item(34, 41).ClearContents()

# Row 35
$ws.Cells.Item(35, 1).Value = 111682658
$ws.Cells.Item(35, 2).Value = 90689
$ws.Cells.Item(35, 5).Value = 5966
$ws.Cells.Item(35, 6).Value = 'Motaggsvamp'
$ws.Cells.Item(35, 7).Value = 'Sarcodon squamosus'
$ws.Cells.Item(35, 8).Value = '(Schaeff.) Quél.'
